# "Add files via upload" - point the sheet's links at the new lxbc.online
# domain instead of the old netlify/google-forms addresses, and leave the
# cursor resting on the last data row (E7) instead of H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the displayed/stored text of the three link cells -------------
# G3: "NỘP ĐƠN GIA NHẬP LXBC" signup form link
$ws.Range("G3").Value = "https://forms.gle/gtgp2jWGgHS126Se8"

# G5: "LXBC Shop" link
$ws.Range("G5").Value = "https://lxbc.online/shop"

# G6: "TKB-LXBC" (schedule) link
$ws.Range("G6").Value = "https://lxbc.online/tkb"

# --- Move the active selection to E7 (also drops the stale topLeftCell) ---
[void]$ws.Range("E7").Select()
